$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill every empty cell in the schedule grid (rows 9-28, columns A-K) with "free"
for ($r = 9; $r -le 28; $r++) {
    for ($c = 1; $c -le 11; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -eq $null) {
            $cell.Value = "free"
        }
    }
}

$ws.Range("N14").Select() | Out-Null
